{"js": "const replacements = [\n  [\"2025-11-18 Tuesday\", \"2025-11-19 Wednesday\"],\n  [\"236\u00f75=47, 1\", \"856\u00f78=107, 0\"],\n  [\"537\u00f74=134, 1\", \"673\u00f78=84, 1\"],\n  [\"413\u00f77=59, 0\", \"142\u00f74=35, 2\"],\n  [\"238\u00f78=29, 6\", \"821\u00f75=164, 1\"],\n  [\"968\u00f72=484, 0\", \"719\u00f73=239, 2\"],\n  [\"621\u00f78=77, 5\", \"655\u00f79=72, 7\"],\n  [\"154\u00f74=38, 2\", \"947\u00f73=315, 2\"],\n  [\"376\u00f78=47, 0\", \"680\u00f75=136, 0\"],\n  [\"401\u00f77=57, 2\", \"603\u00f76=100, 3\"],\n  [\"615\u00f73=205, 0\", \"628\u00f73=209, 1\"],\n  [\"347\u00f72=173, 1\", \"935\u00f76=155, 5\"],\n  [\"172\u00f73=57, 1\", \"161\u00f79=17, 8\"],\n  [\"343\u00f79=38, 1\", \"510\u00f79=56, 6\"],\n  [\"695\u00f79=77, 2\", \"765\u00f75=153, 0\"],\n  [\"156\u00f75=31, 1\", \"318\u00f76=53, 0\"],\n  [\"888\u00f76=148, 0\", \"948\u00f72=474, 0\"],\n  [\"182\u00f78=22, 6\", \"270\u00f76=45, 0\"],\n  [\"673\u00f77=96, 1\", \"194\u00f73=64, 2\"],\n  [\"106\u00f78=13, 2\", \"731\u00f75=146, 1\"],\n  [\"975\u00f79=108, 3\", \"625\u00f78=78, 1\"],\n  [\"406\u00f73=135, 1\", \"812\u00f72=406, 0\"],\n  [\"420\u00f76=70, 0\", \"389\u00f76=64, 5\"],\n  [\"408\u00f76=68, 0\", \"394\u00f75=78, 4\"],\n  [\"114\u00f72=57, 0\", \"729\u00f79=81, 0\"],\n  [\"193\u00f76=32, 1\", \"899\u00f74=224, 3\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$olds = @(\n    \"2025-11-18 Tuesday\"\n    \"236\u00f75=47, 1\"\n    \"537\u00f74=134, 1\"\n    \"413\u00f77=59, 0\"\n    \"238\u00f78=29, 6\"\n    \"968\u00f72=484, 0\"\n    \"621\u00f78=77, 5\"\n    \"154\u00f74=38, 2\"\n    \"376\u00f78=47, 0\"\n    \"401\u00f77=57, 2\"\n    \"615\u00f73=205, 0\"\n    \"347\u00f72=173, 1\"\n    \"172\u00f73=57, 1\"\n    \"343\u00f79=38, 1\"\n    \"695\u00f79=77, 2\"\n    \"156\u00f75=31, 1\"\n    \"888\u00f76=148, 0\"\n    \"182\u00f78=22, 6\"\n    \"673\u00f77=96, 1\"\n    \"106\u00f78=13, 2\"\n    \"975\u00f79=108, 3\"\n    \"406\u00f73=135, 1\"\n    \"420\u00f76=70, 0\"\n    \"408\u00f76=68, 0\"\n    \"114\u00f72=57, 0\"\n    \"193\u00f76=32, 1\"\n)\n\n$news = @(\n    \"2025-11-19 Wednesday\"\n    \"856\u00f78=107, 0\"\n    \"673\u00f78=84, 1\"\n    \"142\u00f74=35, 2\"\n    \"821\u00f75=164, 1\"\n    \"719\u00f73=239, 2\"\n    \"655\u00f79=72, 7\"\n    \"947\u00f73=315, 2\"\n    \"680\u00f75=136, 0\"\n    \"603\u00f76=100, 3\"\n    \"628\u00f73=209, 1\"\n    \"935\u00f76=155, 5\"\n    \"161\u00f79=17, 8\"\n    \"510\u00f79=56, 6\"\n    \"765\u00f75=153, 0\"\n    \"318\u00f76=53, 0\"\n    \"948\u00f72=474, 0\"\n    \"270\u00f76=45, 0\"\n    \"194\u00f73=64, 2\"\n    \"731\u00f75=146, 1\"\n    \"625\u00f78=78, 1\"\n    \"812\u00f72=406, 0\"\n    \"389\u00f76=64, 5\"\n    \"394\u00f75=78, 4\"\n    \"729\u00f79=81, 0\"\n    \"899\u00f74=224, 3\"\n)\n\nfor ($i = 0; $i -lt $olds.Count; $i++) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($olds[$i], $false, $false, $false, $false, $false, $true, 1, $false, $news[$i], 2) | Out-Null\n}"}
